# The workbook currently has a single worksheet "ScenarioA" that holds the
# "Power - Storage" table (headers, units, descriptions and data rows),
# together with its cell comments (A3, B3:B7) and legacy VML comment
# drawing.
#
# The authored change ("Fix exampleStochastic to actually have two
# Scenarios") adds a second, identical scenario sheet named "ScenarioB"
# right after "ScenarioA" - same columns, same data, same comments - so
# that the stochastic example actually has two scenarios to pick from.
#
# Copying the existing sheet duplicates the cell values/styles, the
# column widths, the frozen panes/view, and the cell comments (and their
# VML drawing) in one shot, which is exactly what's needed here.

$wb = $excel.ActiveWorkbook
$scenarioA = $wb.Worksheets.Item("ScenarioA")

# Insert the copy right after ScenarioA.
$scenarioA.Copy($null, $scenarioA)

# The copy is created as "ScenarioA (2)" - rename it to "ScenarioB".
$scenarioB = $wb.Worksheets.Item(2)
$scenarioB.Name = "ScenarioB"

# Copying a sheet makes the new copy the active one; restore ScenarioA as
# the active/selected tab to match the original workbook view.
$scenarioA.Activate()
